$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the first four data rows (Cutoff 0-3), shifting the rest up.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Column A ("Cutoff") is a 0-based sequential index that is
    # regenerated after the shift, so it still reads 0..14 top to bottom.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
